# Update "想去人数" (F column) values on the 展览, 演出, and 全部类型 sheets
# to match the freshly generated data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 301
$ws1.Range("F3").Value = 95
$ws1.Range("F4").Value = 1190
$ws1.Range("F5").Value = 603

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 10

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 301
$ws4.Range("F3").Value = 95
$ws4.Range("F4").Value = 1190
$ws4.Range("F5").Value = 10
$ws4.Range("F6").Value = 603
